$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RH")
$ws.Protect($null, $true, $true, $true)
Write-Output "A done"
